$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
# Row 132 (ALC)
$ws_ALC.Range("H132").Value = 1211.0605
$ws_ALC.Range("I132").Value = 1155.1562
$ws_ALC.Range("J132").Value = 3000
$ws_ALC.Range("K132").Value = 3465.4686
$ws_ALC.Range("L132").Value = 9000
$ws_ALC.Range("M132").Value = -935.4685999999997
$ws_ALC.Range("N132").Value = -14060

$ws_ARM = $wb.Worksheets.Item("ARM")
# Row 32 (ARM)
$ws_ARM.Range("H32").Value = 38370.168
$ws_ARM.Range("I32").Value = 43349.652
$ws_ARM.Range("J32").Value = 6003.5
$ws_ARM.Range("K32").Value = 43349.652
$ws_ARM.Range("L32").Value = 6003.5
$ws_ARM.Range("M32").Value = -43062.652
$ws_ARM.Range("N32").Value = -6577.5

# Row 61 (ARM)
$ws_ARM.Range("H61").Value = 5905.154
$ws_ARM.Range("I61").Value = 2295.5
$ws_ARM.Range("J61").Value = 15702.786
$ws_ARM.Range("K61").Value = 2295.5
$ws_ARM.Range("L61").Value = 15702.786
$ws_ARM.Range("M61").Value = -2083.5
$ws_ARM.Range("N61").Value = -16126.786

# Row 74 (ARM)
$ws_ARM.Range("H74").Value = 5216.067
$ws_ARM.Range("I74").Value = 1811.1666
$ws_ARM.Range("J74").Value = 18835.666
$ws_ARM.Range("K74").Value = 1811.1666
$ws_ARM.Range("L74").Value = 18835.666
$ws_ARM.Range("M74").Value = -937.1666
$ws_ARM.Range("N74").Value = -20583.666

# Row 77 (ARM)
$ws_ARM.Range("H77").Value = 5216.067
$ws_ARM.Range("I77").Value = 1811.1666
$ws_ARM.Range("J77").Value = 18835.666
$ws_ARM.Range("K77").Value = 9055.833000000001
$ws_ARM.Range("L77").Value = 94178.33
$ws_ARM.Range("M77").Value = -4687.833000000001
$ws_ARM.Range("N77").Value = -102914.33

# Row 88 (ARM)
$ws_ARM.Range("H88").Value = 5596.9375
$ws_ARM.Range("I88").Value = 10033.333
$ws_ARM.Range("J88").Value = 2935.1
$ws_ARM.Range("K88").Value = 10033.333
$ws_ARM.Range("L88").Value = 2935.1
$ws_ARM.Range("M88").Value = -9627.333000000001
$ws_ARM.Range("N88").Value = -3747.1

# Row 91 (ARM)
$ws_ARM.Range("H91").Value = 5596.9375
$ws_ARM.Range("I91").Value = 10033.333
$ws_ARM.Range("J91").Value = 2935.1
$ws_ARM.Range("K91").Value = 10033.333
$ws_ARM.Range("L91").Value = 2935.1
$ws_ARM.Range("M91").Value = -8629.333000000001
$ws_ARM.Range("N91").Value = -5743.1

# Row 101 (ARM)
$ws_ARM.Range("H101").Value = 0
$ws_ARM.Range("J101").Value = 0
$ws_ARM.Range("L101").ClearContents()
$ws_ARM.Range("N101").Value = 0

# Row 102 (ARM)
$ws_ARM.Range("H102").Value = 4400
$ws_ARM.Range("I102").Value = 0
$ws_ARM.Range("K102").Value = 0
$ws_ARM.Range("M102").ClearContents()

# Row 136 (ARM)
$ws_ARM.Range("H136").Value = 5905.154
$ws_ARM.Range("I136").Value = 2295.5
$ws_ARM.Range("J136").Value = 15702.786
$ws_ARM.Range("K136").Value = 6886.5
$ws_ARM.Range("L136").Value = 47108.358
$ws_ARM.Range("M136").Value = -4336.5
$ws_ARM.Range("N136").Value = -52208.358

$ws_BSM = $wb.Worksheets.Item("BSM")
# Row 86 (BSM)
$ws_BSM.Range("H86").Value = 1927.4717
$ws_BSM.Range("I86").Value = 1789.4894
$ws_BSM.Range("K86").Value = 1789.4894
$ws_BSM.Range("M86").Value = -666.4893999999999

# Row 89 (BSM)
$ws_BSM.Range("H89").Value = 1927.4717
$ws_BSM.Range("I89").Value = 1789.4894
$ws_BSM.Range("K89").Value = 8947.447
$ws_BSM.Range("M89").Value = -3331.447

# Row 103 (BSM)
$ws_BSM.Range("H103").Value = 17323.666
$ws_BSM.Range("J103").Value = 17323.666
$ws_BSM.Range("L103").Value = 17323.666
$ws_BSM.Range("N103").Value = -19667.666

# Row 105 (BSM)
$ws_BSM.Range("H105").Value = 5679.7856
$ws_BSM.Range("I105").Value = 5592.4546
$ws_BSM.Range("J105").Value = 6000
$ws_BSM.Range("K105").Value = 5592.4546
$ws_BSM.Range("L105").Value = 6000
$ws_BSM.Range("M105").Value = -3845.4546
$ws_BSM.Range("N105").Value = -9494

$ws_CUL = $wb.Worksheets.Item("CUL")
# Row 39 (CUL)
$ws_CUL.Range("H39").Value = 4400
$ws_CUL.Range("J39").Value = 4400
$ws_CUL.Range("L39").Value = 13200
$ws_CUL.Range("N39").Value = -13788

# Row 75 (CUL)
$ws_CUL.Range("H75").Value = 4091.3
$ws_CUL.Range("I75").Value = 1971
$ws_CUL.Range("J75").Value = 5000
$ws_CUL.Range("K75").Value = 5913
$ws_CUL.Range("L75").Value = 15000
$ws_CUL.Range("M75").Value = -4915
$ws_CUL.Range("N75").Value = -16996

# Row 78 (CUL)
$ws_CUL.Range("H78").Value = 4091.3
$ws_CUL.Range("I78").Value = 1971
$ws_CUL.Range("J78").Value = 5000
$ws_CUL.Range("K78").Value = 17739
$ws_CUL.Range("L78").Value = 45000
$ws_CUL.Range("M78").Value = -12747
$ws_CUL.Range("N78").Value = -54984

# Row 92 (CUL)
$ws_CUL.Range("H92").Value = 249.5
$ws_CUL.Range("I92").Value = 249.5
$ws_CUL.Range("J92").Value = 0
$ws_CUL.Range("K92").Value = 748.5
$ws_CUL.Range("L92").Value = 0
$ws_CUL.Range("M92").ClearContents()
$ws_CUL.Range("N92").Value = 499.5

# Row 102 (CUL)
$ws_CUL.Range("H102").Value = 27513
$ws_CUL.Range("J102").Value = 50000
$ws_CUL.Range("L102").Value = 150000
$ws_CUL.Range("N102").Value = -154868

# Row 110 (CUL)
$ws_CUL.Range("H110").Value = 3284.5
$ws_CUL.Range("J110").Value = 0
$ws_CUL.Range("L110").Value = 0
$ws_CUL.Range("N110").ClearContents()

# Row 119 (CUL)
$ws_CUL.Range("H119").Value = 6305.8
$ws_CUL.Range("I119").Value = 529
$ws_CUL.Range("J119").Value = 7750
$ws_CUL.Range("K119").Value = 1587
$ws_CUL.Range("L119").Value = 23250
$ws_CUL.Range("M119").Value = 3251
$ws_CUL.Range("N119").Value = -32926

$ws_GSM = $wb.Worksheets.Item("GSM")
# Row 132 (GSM)
$ws_GSM.Range("H132").Value = 49264.5
$ws_GSM.Range("I132").Value = 85903.164
$ws_GSM.Range("J132").Value = 12625.833
$ws_GSM.Range("K132").Value = 257709.492
$ws_GSM.Range("L132").Value = 37877.499
$ws_GSM.Range("M132").Value = -255179.492
$ws_GSM.Range("N132").Value = -42937.499

$ws_LTW = $wb.Worksheets.Item("LTW")
# Row 68 (LTW)
$ws_LTW.Range("H68").Value = 3729.6538
$ws_LTW.Range("I68").Value = 2840.8333
$ws_LTW.Range("J68").Value = 4491.5
$ws_LTW.Range("K68").Value = 2840.8333
$ws_LTW.Range("L68").Value = 4491.5
$ws_LTW.Range("M68").Value = -2091.8333
$ws_LTW.Range("N68").Value = -5989.5

# Row 71 (LTW)
$ws_LTW.Range("H71").Value = 3729.6538
$ws_LTW.Range("I71").Value = 2840.8333
$ws_LTW.Range("J71").Value = 4491.5
$ws_LTW.Range("K71").Value = 14204.1665
$ws_LTW.Range("L71").Value = 22457.5
$ws_LTW.Range("M71").Value = -10460.1665
$ws_LTW.Range("N71").Value = -29945.5

# Row 100 (LTW)
$ws_LTW.Range("H100").Value = 4933.3335
$ws_LTW.Range("I100").Value = 4425
$ws_LTW.Range("J100").Value = 9000
$ws_LTW.Range("K100").Value = 4425
$ws_LTW.Range("L100").Value = 9000
$ws_LTW.Range("M100").Value = -3884
$ws_LTW.Range("N100").Value = -10082

$ws_WVR = $wb.Worksheets.Item("WVR")
# Row 81 (WVR)
$ws_WVR.Range("H81").Value = 25004064
$ws_WVR.Range("J81").Value = 33337502
$ws_WVR.Range("L81").Value = 66675004
$ws_WVR.Range("N81").Value = -66677126

# Row 84 (WVR)
$ws_WVR.Range("H84").Value = 25004064
$ws_WVR.Range("J84").Value = 33337502
$ws_WVR.Range("L84").Value = 333375020
$ws_WVR.Range("N84").Value = -333385628
